# Applies the cell-level updates described by the commit diff.
# The workbook stores numeric-looking "Price" values (column D) as literal
# text strings (inlineStr in the original OOXML), so for those cells we
# force a text number format before assignment (and restore the default
# "Normal" style afterwards) to avoid Excel silently re-typing them as
# numbers, which would change both the stored type and the style index.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextNumberCell($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

function Set-TextCell($cell, $value) {
    $ws.Range($cell).Value = $value
}

Set-TextNumberCell "D2" "251.03"
Set-TextNumberCell "D3" "22.77"
Set-TextNumberCell "D4" "5.426"
Set-TextNumberCell "D5" "0.05676"
Set-TextNumberCell "D6" "3.411"
Set-TextNumberCell "D7" "6.369"
Set-TextNumberCell "D8" "0.8131"
Set-TextNumberCell "D9" "0.9301"
Set-TextCell "B10" "One"
Set-TextCell "C10" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextNumberCell "D10" "0.01158"
Set-TextCell "E10" "9OneONEBestin24h"
Set-TextCell "B11" "WazirX"
Set-TextCell "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextNumberCell "D11" "0.1441"
Set-TextCell "E11" "10WazirXWRX"
Set-TextCell "B12" "MandalaExchangeToken"
Set-TextCell "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextNumberCell "D12" "0.07459"
Set-TextCell "E12" "11MandalaExchangeTokenMDX"
Set-TextCell "B13" "LiechtensteinCryptoassetsExchange"
Set-TextCell "C13" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextNumberCell "D13" "0.03163"
Set-TextCell "E13" "12LiechtensteinCryptoassetsExchangeLCX"
Set-TextCell "B14" "BitrueCoin"
Set-TextCell "C14" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextNumberCell "D14" "0.03074"
Set-TextCell "E14" "13BitrueCoinBTR"
Set-TextCell "B15" "BitMartToken"
Set-TextCell "C15" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextNumberCell "D15" "0.09353"
Set-TextCell "E15" "14BitMartTokenBMX"
Set-TextCell "B16" "MCDex"
Set-TextCell "C16" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextNumberCell "D16" "3.729"
Set-TextCell "E16" "15MCDexMCB"
Set-TextCell "B17" "BitForexToken"
Set-TextCell "C17" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextNumberCell "D17" "0.001585"
Set-TextCell "E17" "16BitForexTokenBF"
Set-TextCell "B18" "CoinExToken"
Set-TextCell "C18" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextNumberCell "D18" "0.04775"
Set-TextCell "E18" "17CoinExTokenCET"
Set-TextNumberCell "D20" "0.005044"
Set-TextNumberCell "D21" "0.001030"
Set-TextNumberCell "D23" "3.713"
Set-TextNumberCell "D24" "2.170"
Set-TextNumberCell "D28" "0.0002999"
Set-TextNumberCell "D40" "0.04022"
Set-TextNumberCell "D41" "0.006788"
Set-TextNumberCell "D42" "0.1070"
Set-TextNumberCell "D43" "0.002709"
Set-TextNumberCell "D44" "0.007552"
Set-TextNumberCell "D45" "0.00005759"
Set-TextNumberCell "D47" "0.4999"
Set-TextCell "E48" "47BOLOBOLO"
